$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5105
$ws.Range("I3").Value = 5324
$ws.Range("E4").Value = 1966
$ws.Range("F4").Value = 1866
$ws.Range("H4").Value = 1673
$ws.Range("I4").Value = 1221
$ws.Range("I5").Value = 495
$ws.Range("I6").Value = 5811
$ws.Range("E7").Value = 25970
$ws.Range("F7").Value = 24055
$ws.Range("H7").Value = 25983
$ws.Range("I7").Value = 17956

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 54
$ws.Range("I7").Value = 191

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I2").Value = 39
$ws.Range("I7").Value = 106

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 57

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 180
$ws.Range("I3").Value = 187
$ws.Range("I7").Value = 575

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 120
$ws.Range("I7").Value = 335

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 256
$ws.Range("I4").Value = 44
$ws.Range("I5").Value = 20
$ws.Range("I7").Value = 706

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I3").Value = 49
$ws.Range("I7").Value = 150

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 134
$ws.Range("I6").Value = 117
$ws.Range("I7").Value = 402

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 143
$ws.Range("I4").Value = 71
$ws.Range("I6").Value = 123
$ws.Range("I7").Value = 580
$ws.Range("I8").Value = 1083
$ws.Range("I9").Value = 79
$ws.Range("I11").Value = 268
$ws.Range("I14").Value = 106
$ws.Range("I19").Value = 494
$ws.Range("I22").Value = 46
$ws.Range("I23").Value = 176
$ws.Range("I27").Value = 168
$ws.Range("I29").Value = 1141
$ws.Range("I30").Value = 57
$ws.Range("I33").Value = 821
$ws.Range("I36").Value = 236
$ws.Range("I37").Value = 575
$ws.Range("I42").Value = 606
$ws.Range("I43").Value = 142
$ws.Range("I44").Value = 128
$ws.Range("I46").Value = 37
$ws.Range("I47").Value = 120
$ws.Range("I48").Value = 246
$ws.Range("I51").Value = 195
$ws.Range("I52").Value = 399
$ws.Range("I53").Value = 185
$ws.Range("I54").Value = 385
$ws.Range("E63").Value = 313
$ws.Range("F63").Value = 157
$ws.Range("H63").Value = 220
$ws.Range("I63").Value = 65
$ws.Range("I65").Value = 402
$ws.Range("I67").Value = 706
$ws.Range("I72").Value = 66
$ws.Range("I73").Value = 156
$ws.Range("I76").Value = 268
$ws.Range("I78").Value = 255
$ws.Range("I79").Value = 515
$ws.Range("I80").Value = 61
$ws.Range("I83").Value = 378
$ws.Range("I84").Value = 150
$ws.Range("I85").Value = 811
$ws.Range("I86").Value = 109
$ws.Range("I88").Value = 166
$ws.Range("I94").Value = 180
$ws.Range("I95").Value = 290
$ws.Range("I96").Value = 191
$ws.Range("I97").Value = 139
$ws.Range("I99").Value = 335
$ws.Range("E101").Value = 25970
$ws.Range("F101").Value = 24055
$ws.Range("H101").Value = 25983
$ws.Range("I101").Value = 17956

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I3").Value = 143
$ws.Range("I7").Value = 378

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 100
$ws.Range("I7").Value = 290

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 307
$ws.Range("I6").Value = 257
$ws.Range("I7").Value = 821

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 184
$ws.Range("I7").Value = 385

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 340
$ws.Range("I3").Value = 395
$ws.Range("I4").Value = 59
$ws.Range("I6").Value = 306
$ws.Range("I7").Value = 1141

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 149
$ws.Range("I5").Value = 12
$ws.Range("I6").Value = 135
$ws.Range("I7").Value = 494

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I6").Value = 40
$ws.Range("I7").Value = 128

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I3").Value = 49
$ws.Range("I7").Value = 246

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I6").Value = 121
$ws.Range("I7").Value = 268

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 222
$ws.Range("I3").Value = 323
$ws.Range("I7").Value = 811

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 157
$ws.Range("I7").Value = 606

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I3").Value = 64
$ws.Range("I6").Value = 92
$ws.Range("I7").Value = 255

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("I6").Value = 14
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 62
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 147
$ws.Range("I3").Value = 166
$ws.Range("I5").Value = 16
$ws.Range("I6").Value = 152
$ws.Range("I7").Value = 515

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I4").Value = 31
$ws.Range("I6").Value = 135

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 75
$ws.Range("I7").Value = 236

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 108
$ws.Range("I3").Value = 146
$ws.Range("I7").Value = 399

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 180

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I3").Value = 36
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I3").Value = 51
$ws.Range("I7").Value = 268

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I2").Value = 31
$ws.Range("I3").Value = 27
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 156

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I2").Value = 47
$ws.Range("I7").Value = 143

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I6").Value = 85
$ws.Range("I7").Value = 139

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 51
$ws.Range("I7").Value = 166

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 340
$ws.Range("I4").Value = 62
$ws.Range("I6").Value = 346
$ws.Range("I7").Value = 1083

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I3").Value = 32
$ws.Range("I6").Value = 66
$ws.Range("I7").Value = 168

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I4").Value = 54
$ws.Range("I7").Value = 109

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I4").Value = 22
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I3").Value = 25
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 41
$ws.Range("I7").Value = 185

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I2").Value = 17
$ws.Range("I7").Value = 46

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 66

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 191
$ws.Range("I6").Value = 152
$ws.Range("I7").Value = 580

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I3").Value = 19
$ws.Range("I7").Value = 71
